# Apply updated dSF (column F) values per row, as recorded in the commit
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -12
    3  = -10
    6  = -5
    9  = -5
    10 = -11
    12 = -4
    13 = -8
    14 = 4
    15 = -14
    17 = 0
    22 = -9
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
